$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:I1 - copy style from existing header cell E1
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:I1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("F1").Value = "Average EK between LR and HR"
$ws.Range("G1").Value = "Average EK between Model Output and HR"
$ws.Range("H1").Value = "Average Acc between LR and HR"
$ws.Range("I1").Value = "Average Acc between Model Output and HR"

# Data values for rows 2-7, columns F, G, H, I
$data = @{
    2 = @(8.522042723894119, 3.53161281367143, 0.9627068261299901, 0.9948936274063523)
    3 = @(20.49433722813924, 11.80521180232366, 0.7566254261269123, 0.9395661467287658)
    4 = @(32.23843422253927, 23.81968121210734, 0.4969428583524946, 0.7006880093674811)
    5 = @(8.559897836446762, 3.493587876955668, 0.961297985508244, 0.9949249328997131)
    6 = @(20.58161719322204, 11.85936235427856, 0.7489690695827318, 0.9372976741223997)
    7 = @(32.14891414324443, 23.48789763450623, 0.4872306777600225, 0.7050356599753942)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("F$row").Value = $vals[0]
    $ws.Range("G$row").Value = $vals[1]
    $ws.Range("H$row").Value = $vals[2]
    $ws.Range("I$row").Value = $vals[3]
}
